$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need an explicit Text
# number format applied first, otherwise Excel auto-converts the assigned
# string into a numeric value (losing e.g. trailing zeros).
$textCells = @("D5", "D6", "D12", "D15", "D21", "D22", "D24", "D30", "D33", "D34", "D36", "D37", "D40", "D41", "D44", "D45", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "58.507.81"
$ws.Range("E2").Value = "  -3.92%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "533.80"
$ws.Range("D6").Value = "134.63"
$ws.Range("E6").Value = "  -8.48%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "3.192.43"
$ws.Range("E8").Value = "  -4.89%  "
$ws.Range("E9").Value = "  -4.99%  "
$ws.Range("E10").Value = "  -6.53%  "
$ws.Range("E11").Value = "  -7.14%  "
$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  -5.11%  "
$ws.Range("E13").Value = "  -4.91%  "
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Value = "25.70"
$ws.Range("E15").Value = "  -7.99%  "
$ws.Range("E16").Value = "  -4.90%  "
$ws.Range("D17").Value = "58.602.91"
$ws.Range("E17").Value = "  -3.89%  "
$ws.Range("E18").Value = "  -7.72%  "
$ws.Range("E19").Value = "  -6.54%  "
$ws.Range("E20").Value = "  -8.60%  "
$ws.Range("D21").Value = "8.12"
$ws.Range("E21").Value = "  -8.89%  "
$ws.Range("D22").Value = "358.65"
$ws.Range("E22").Value = "  -4.69%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "69.66"
$ws.Range("E24").Value = "  -7.05%  "
$ws.Range("E25").Value = "  -7.82%  "
$ws.Range("E26").Value = "  -5.07%  "
$ws.Range("E27").Value = "  -3.41%  "
$ws.Range("D28").Value = "0.0₃0950"
$ws.Range("E28").Value = "  -12.01%  "
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").Value = "7.03"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  -8.31%  "
$ws.Range("D33").Value = "7.00"
$ws.Range("E33").Value = "  -8.96%  "
$ws.Range("D34").Value = "21.65"
$ws.Range("E34").Value = "  -5.21%  "
$ws.Range("E35").Value = "  -9.10%  "
$ws.Range("D36").Value = "4.95"
$ws.Range("E36").Value = "  -7.35%  "
$ws.Range("D37").Value = "160.85"
$ws.Range("E37").Value = "  -4.97%  "
$ws.Range("E38").Value = "  -7.04%  "
$ws.Range("E39").Value = "  -8.49%  "
$ws.Range("D40").Value = "25.82"
$ws.Range("E40").Value = "  -10.77%  "
$ws.Range("D41").Value = "0.0704"
$ws.Range("E41").Value = "  -6.15%  "
$ws.Range("E42").Value = "  -5.04%  "
$ws.Range("E43").Value = "  -3.66%  "
$ws.Range("D44").Value = "0.709"
$ws.Range("E44").Value = "  -6.71%  "
$ws.Range("D45").Value = "1.09"
$ws.Range("E45").Value = "  -3.72%  "
$ws.Range("E46").Value = "  -6.85%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  -7.84%  "
$ws.Range("D49").Value = "2.278.71"
$ws.Range("E49").Value = "  -8.72%  "
$ws.Range("E50").Value = "  -6.22%  "
$ws.Range("D51").Value = "20.44"
$ws.Range("E51").Value = "  -9.68%  "
